$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4353, 4481, 4524, 4524, 4524, 4850, 4850, 4850, 4850, 4884, 4884)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
